# Insert a new weekly "Perejil" (parsley) price record for the
# Vega Central Mapocho de Santiago market ahead of the row currently
# holding 2022-04-25 (serial 44664), pushing the existing rows 361-379
# down to 362-380.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 361 (shifts rows 361:379 -> 362:380,
# carrying their formatting/styles with them).
$ws.Rows.Item(361).Insert()

# Populate the new row with the new record's data.
$ws.Cells.Item(361, 1).Value2 = 9
$ws.Cells.Item(361, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(361, 3).Value2 = "Metropolitana"
$ws.Cells.Item(361, 4).Value2 = 44706
$ws.Cells.Item(361, 5).Value2 = 13
$ws.Cells.Item(361, 6).Value2 = 100112044
$ws.Cells.Item(361, 7).Value2 = "Perejil"
$ws.Cells.Item(361, 8).Value2 = "Sin especificar"
$ws.Cells.Item(361, 9).Value2 = "Primera"
$ws.Cells.Item(361, 10).Value2 = 62
$ws.Cells.Item(361, 11).Value2 = 8000
$ws.Cells.Item(361, 12).Value2 = 9000
$ws.Cells.Item(361, 13).Value2 = 8516
$ws.Cells.Item(361, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(361, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(361, 16).Value2 = 2839
$ws.Cells.Item(361, 17).Value2 = 3
$ws.Cells.Item(361, 18).Value2 = "Hortaliza"
